# Dicionario de dados - "arrumando pequeno erro no dicionario de dados"
#
# 1) The "Tamanho" (size) cell for the idUser field (D4) used to contain the
#    raw numeric literal 4294967295 (the max value of an unsigned 4-byte
#    integer). That was the "small error": it should instead simply describe
#    the storage size, like the other numeric fields in the sheet, so it is
#    replaced with the text "4 Bytes".
# 2) The saved selection/cursor position in the sheet moves from M8 to K6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "4 Bytes"

$null = $ws.Range("K6").Select()
